$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "2019-IKA/WQI"
$ws.Range("F1").Value = "2019-category"
$ws.Range("G1").Value = "2019-class"
$ws.Range("H1").Value = "2020-IKA/WQI"
$ws.Range("I1").Value = "2020-category"
$ws.Range("J1").Value = "2020-class"
